# Update "想去人数" (interest count, column F) figures that changed between
# scrapes of the Bilibili "show" listings, as published to gh-pages.
#
# Sheet "展览"   (Exhibitions)
# Sheet "演出"   (Performances)
# Sheet "全部类型" (All types)
# Sheet "本地生活" (Local life) is unaffected.

$wb = $excel.ActiveWorkbook

function Set-F {
    param($sheet, $row, $value)
    $sheet.Cells.Item($row, 6).Value = $value
}

# --- 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
Set-F $ws1 3  150
Set-F $ws1 7  4209
Set-F $ws1 8  345
Set-F $ws1 9  220
Set-F $ws1 13 8
Set-F $ws1 17 1516
Set-F $ws1 18 1396
Set-F $ws1 23 427
Set-F $ws1 24 87

# --- 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
Set-F $ws2 7 5

# --- 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
Set-F $ws4 4  150
Set-F $ws4 9  4209
Set-F $ws4 10 345
Set-F $ws4 11 220
Set-F $ws4 19 8
Set-F $ws4 22 5
Set-F $ws4 24 1516
Set-F $ws4 25 1396
Set-F $ws4 31 427
Set-F $ws4 32 87
